$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 6336.1113
$ws.Range("I76").Value = 5050
$ws.Range("J76").Value = 7365
$ws.Range("K76").Value = 5050
$ws.Range("L76").Value = 7365
$ws.Range("M76").Value = -4735
$ws.Range("N76").Value = -7995
$ws.Range("H79").Value = 6336.1113
$ws.Range("I79").Value = 5050
$ws.Range("J79").Value = 7365
$ws.Range("K79").Value = 5050
$ws.Range("L79").Value = 7365
$ws.Range("M79").Value = -3958
$ws.Range("N79").Value = -9549
$ws.Range("H104").Value = 1425
$ws.Range("I104").Value = 1850
$ws.Range("J104").Value = 150
$ws.Range("K104").Value = 5550
$ws.Range("L104").Value = 450
$ws.Range("M104").Value = -3803
$ws.Range("N104").Value = -3944
$ws.Range("H138").Value = 2634.3457
$ws.Range("I138").Value = 3249.6667
$ws.Range("K138").Value = 9749.000100000001
$ws.Range("M138").Value = -4609.000100000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 5313.1
$ws.Range("J46").Value = 5791.375
$ws.Range("L46").Value = 5791.375
$ws.Range("N46").Value = -6429.375
$ws.Range("H61").Value = 57444.453
$ws.Range("I61").Value = 2082.4075
$ws.Range("J61").Value = 157096.14
$ws.Range("K61").Value = 2082.4075
$ws.Range("L61").Value = 157096.14
$ws.Range("M61").Value = -1870.4075
$ws.Range("N61").Value = -157520.14
$ws.Range("H74").Value = 12068.564
$ws.Range("J74").Value = 46880.11
$ws.Range("L74").Value = 46880.11
$ws.Range("N74").Value = -48628.11
$ws.Range("H77").Value = 12068.564
$ws.Range("J77").Value = 46880.11
$ws.Range("L77").Value = 234400.55
$ws.Range("N77").Value = -243136.55
$ws.Range("H97").Value = 433.9655
$ws.Range("I97").Value = 467.72223
$ws.Range("J97").Value = 378.72726
$ws.Range("K97").Value = 467.72223
$ws.Range("L97").Value = 378.72726
$ws.Range("M97").Value = 28.27776999999998
$ws.Range("N97").Value = -1370.72726
$ws.Range("H102").Value = 2738.6155
$ws.Range("I102").Value = 1716.8334
$ws.Range("J102").Value = 15000
$ws.Range("K102").Value = 1716.8334
$ws.Range("L102").Value = 15000
$ws.Range("M102").Value = -94.83339999999998
$ws.Range("N102").Value = -18244
$ws.Range("H122").Value = 1340480.4
$ws.Range("I122").Value = 2263851
$ws.Range("K122").Value = 6791553
$ws.Range("M122").Value = -6789103
$ws.Range("H132").Value = 3861645
$ws.Range("I132").Value = 4677.2
$ws.Range("K132").Value = 14031.6
$ws.Range("M132").Value = -11501.6
$ws.Range("H133").Value = 49947.75
$ws.Range("J133").Value = 49947.75
$ws.Range("L133").Value = 49947.75
$ws.Range("N133").Value = -55007.75
$ws.Range("H136").Value = 57444.453
$ws.Range("I136").Value = 2082.4075
$ws.Range("J136").Value = 157096.14
$ws.Range("K136").Value = 6247.2225
$ws.Range("L136").Value = 471288.42
$ws.Range("M136").Value = -3697.2225
$ws.Range("N136").Value = -476388.42

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2615.3447
$ws.Range("I94").Value = 1536.3043
$ws.Range("K94").Value = 1536.3043
$ws.Range("M94").Value = -1085.3043

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 13076.667
$ws.Range("I16").Value = 9384.714
$ws.Range("J16").Value = 25998.5
$ws.Range("K16").Value = 9384.714
$ws.Range("L16").Value = 25998.5
$ws.Range("M16").Value = -9097.714
$ws.Range("N16").Value = -26572.5
$ws.Range("H107").Value = 704.5833
$ws.Range("J107").Value = 416.06897
$ws.Range("L107").Value = 416.06897
$ws.Range("N107").Value = -4256.06897
$ws.Range("H113").Value = 13076.667
$ws.Range("I113").Value = 9384.714
$ws.Range("J113").Value = 25998.5
$ws.Range("K113").Value = 9384.714
$ws.Range("L113").Value = 25998.5
$ws.Range("M113").Value = -7214.714
$ws.Range("N113").Value = -30338.5
$ws.Range("H121").Value = 47001
$ws.Range("J121").Value = 47001
$ws.Range("L121").Value = 47001
$ws.Range("N121").Value = -49621
$ws.Range("H122").Value = 2156.5217
$ws.Range("I122").Value = 1820.35
$ws.Range("J122").Value = 4397.6665
$ws.Range("K122").Value = 5461.049999999999
$ws.Range("L122").Value = 13192.9995
$ws.Range("M122").Value = -3011.049999999999
$ws.Range("N122").Value = -18092.9995
$ws.Range("H132").Value = 2749.875
$ws.Range("I132").Value = 2749.875
$ws.Range("K132").Value = 8249.625
$ws.Range("M132").Value = -5719.625
$ws.Range("H134").Value = 29417780
$ws.Range("I134").Value = 2093.5454
$ws.Range("J134").Value = 83346536
$ws.Range("K134").Value = 6280.6362
$ws.Range("L134").Value = 250039608
$ws.Range("M134").Value = -3745.6362
$ws.Range("N134").Value = -250044678

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 4868.66
$ws.Range("I134").Value = 1615.909
$ws.Range("J134").Value = 5786.1025
$ws.Range("K134").Value = 4847.727000000001
$ws.Range("L134").Value = 17358.3075
$ws.Range("M134").Value = 222.2729999999992
$ws.Range("N134").Value = -27498.3075
$ws.Range("H139").Value = 23376.857
$ws.Range("I139").Value = 36534.5
$ws.Range("J139").Value = 5833.3335
$ws.Range("K139").Value = 109603.5
$ws.Range("L139").Value = 17500.0005
$ws.Range("M139").Value = -104463.5
$ws.Range("N139").Value = -27780.0005

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 20655.363
$ws.Range("J70").Value = 19502.25
$ws.Range("L70").Value = 19502.25
$ws.Range("N70").Value = -20042.25
$ws.Range("H73").Value = 20655.363
$ws.Range("J73").Value = 19502.25
$ws.Range("L73").Value = 19502.25
$ws.Range("N73").Value = -21374.25
$ws.Range("H80").Value = 20027.857
$ws.Range("I80").Value = 19549.125
$ws.Range("J80").Value = 20666.166
$ws.Range("K80").Value = 19549.125
$ws.Range("L80").Value = 20666.166
$ws.Range("M80").Value = -18551.125
$ws.Range("N80").Value = -22662.166
$ws.Range("H83").Value = 20027.857
$ws.Range("I83").Value = 19549.125
$ws.Range("J83").Value = 20666.166
$ws.Range("K83").Value = 97745.625
$ws.Range("L83").Value = 103330.83
$ws.Range("M83").Value = -92753.625
$ws.Range("N83").Value = -113314.83
$ws.Range("H102").Value = 4633318
$ws.Range("I102").Value = 9621054
$ws.Range("J102").Value = 1849.2142
$ws.Range("K102").Value = 9621054
$ws.Range("L102").Value = 1849.2142
$ws.Range("M102").Value = -9619432
$ws.Range("N102").Value = -5093.2142
$ws.Range("H107").Value = 713.0833
$ws.Range("J107").Value = 1279.8
$ws.Range("L107").Value = 1279.8
$ws.Range("N107").Value = -5119.8
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H110").Value = 225000
$ws.Range("J110").Value = 225000
$ws.Range("L110").Value = 225000
$ws.Range("N110").Value = -233180
$ws.Range("H132").Value = 11074.571
$ws.Range("I132").Value = 4758.8
$ws.Range("K132").Value = 14276.4
$ws.Range("M132").Value = -11746.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3182.6667
$ws.Range("I82").Value = 3322.2307
$ws.Range("K82").Value = 3322.2307
$ws.Range("M82").Value = -2961.2307
$ws.Range("H85").Value = 3182.6667
$ws.Range("I85").Value = 3322.2307
$ws.Range("K85").Value = 3322.2307
$ws.Range("M85").Value = -2074.2307
$ws.Range("H132").Value = 4021697.5
$ws.Range("I132").Value = 2423.8572
$ws.Range("J132").Value = 13400003
$ws.Range("K132").Value = 7271.571599999999
$ws.Range("L132").Value = 40200009
$ws.Range("M132").Value = -4741.571599999999
$ws.Range("N132").Value = -40205069
$ws.Range("H136").Value = 16184.069
$ws.Range("I136").Value = 17945.385
$ws.Range("J136").Value = 14753
$ws.Range("K136").Value = 53836.155
$ws.Range("L136").Value = 44259
$ws.Range("M136").Value = -51286.155
$ws.Range("N136").Value = -49359

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1714.8572
$ws.Range("I96").Value = 1000
$ws.Range("J96").Value = 1834
$ws.Range("K96").Value = 1000
$ws.Range("L96").Value = 1834
$ws.Range("M96").Value = 373
$ws.Range("N96").Value = -4580
$ws.Range("H107").Value = 1964
$ws.Range("I107").Value = 2208.1667
$ws.Range("K107").Value = 6624.500100000001
$ws.Range("M107").Value = -4704.500100000001
$ws.Range("H122").Value = 399995.53
$ws.Range("I122").Value = 496212.47
$ws.Range("J122").Value = 6380.8184
$ws.Range("K122").Value = 1488637.41
$ws.Range("L122").Value = 19142.4552
$ws.Range("M122").Value = -1486187.41
$ws.Range("N122").Value = -24042.4552
$ws.Range("H132").Value = 17804.03
$ws.Range("I132").Value = 6969.5356
$ws.Range("J132").Value = 61142
$ws.Range("K132").Value = 20908.6068
$ws.Range("L132").Value = 183426
$ws.Range("M132").Value = -18378.6068
$ws.Range("N132").Value = -188486
